$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 — Archer Aviation Inc. (ACHR)
$ws.Range("E2").Value = 45.4
$ws.Range("G2").Value = 40
$ws.Range("K2").Value = 63
$ws.Range("M2").Value = "📈 매수 관찰 구간입니다."
$ws.Range("N2").Value = 85.87127175646313

# Row 3 — Joby Aviation, Inc. (JOBY)
$ws.Range("E3").Value = 46.2
$ws.Range("G3").Value = 30
$ws.Range("K3").Value = 60
$ws.Range("N3").Value = 85.87127175646313
